# Update buddingtonite analysis using XPP matrix correction method
# Rewrites the numeric detection-limit data in rows 2, 5 and 8 on both
# worksheets with the recalculated (XPP matrix correction) values.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("buddingtonite_2_bg_detlim")
$ws2 = $wb.Worksheets.Item("buddingtonite_3_bg_apf_detlim")

# --- Sheet: buddingtonite_2_bg_detlim ---

# Row 2
$ws1.Range("B2").Value = 0.04
$ws1.Range("C2").Value = 0.037
$ws1.Range("D2").Value = 0.04
$ws1.Range("E2").Value = 0.035
$ws1.Range("F2").Value = 0.036
$ws1.Range("G2").Value = 0.038
$ws1.Range("H2").Value = 0.002
$ws1.Range("I2").Value = 0.035
$ws1.Range("J2").Value = 0.04

# Row 5
$ws1.Range("B5").Value = 0.063
$ws1.Range("C5").Value = 0.059
$ws1.Range("D5").Value = 0.062
$ws1.Range("E5").Value = 0.055
$ws1.Range("F5").Value = 0.057
$ws1.Range("G5").Value = 0.059
$ws1.Range("H5").Value = 0.004
$ws1.Range("I5").Value = 0.055
$ws1.Range("J5").Value = 0.063

# Row 8
$ws1.Range("B8").Value = 0.066
$ws1.Range("C8").Value = 0.061
$ws1.Range("D8").Value = 0.065
$ws1.Range("E8").Value = 0.057
$ws1.Range("F8").Value = 0.059
$ws1.Range("G8").Value = 0.062
$ws1.Range("H8").Value = 0.004
$ws1.Range("I8").Value = 0.057
$ws1.Range("J8").Value = 0.066

# --- Sheet: buddingtonite_3_bg_apf_detlim ---

# Row 2
$ws2.Range("B2").Value = 0.049
$ws2.Range("C2").Value = 0.046
$ws2.Range("D2").Value = 0.049
$ws2.Range("E2").Value = 0.043
$ws2.Range("F2").Value = 0.044
$ws2.Range("G2").Value = 0.046
$ws2.Range("H2").Value = 0.003
$ws2.Range("I2").Value = 0.043
$ws2.Range("J2").Value = 0.049

# Row 5
$ws2.Range("B5").Value = 0.077
$ws2.Range("C5").Value = 0.07199999999999999
$ws2.Range("D5").Value = 0.076
$ws2.Range("E5").Value = 0.067
$ws2.Range("F5").Value = 0.07000000000000001
$ws2.Range("G5").Value = 0.07199999999999999
$ws2.Range("H5").Value = 0.004
$ws2.Range("I5").Value = 0.067
$ws2.Range("J5").Value = 0.077

# Row 8
$ws2.Range("B8").Value = 0.08
$ws2.Range("C8").Value = 0.074
$ws2.Range("D8").Value = 0.08
$ws2.Range("E8").Value = 0.07000000000000001
$ws2.Range("F8").Value = 0.07199999999999999
$ws2.Range("G8").Value = 0.075
$ws2.Range("H8").Value = 0.005
$ws2.Range("I8").Value = 0.07000000000000001
$ws2.Range("J8").Value = 0.08
